$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2560.2083
$ws.Range("I106").Value = 2567.5
$ws.Range("J106").Value = 2550
$ws.Range("K106").Value = 2567.5
$ws.Range("L106").Value = 2550
$ws.Range("M106").Value = -1936.5
$ws.Range("N106").Value = -3812

$ws.Range("H107").Value = 1155.4642
$ws.Range("I107").Value = 1439.0714
$ws.Range("K107").Value = 1439.0714
$ws.Range("M107").Value = 480.9286

$ws.Range("H138").Value = 2654.6562
$ws.Range("I138").Value = 2595.9524
$ws.Range("J138").Value = 2766.7273
$ws.Range("K138").Value = 7787.8572
$ws.Range("L138").Value = 8300.1819
$ws.Range("M138").Value = -2647.8572
$ws.Range("N138").Value = -18580.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2102799
$ws.Range("I2").Value = 2332.8572
$ws.Range("K2").Value = 2332.8572
$ws.Range("M2").Value = -2219.8572

$ws.Range("H110").Value = 949.0606
$ws.Range("I110").Value = 851.2727
$ws.Range("J110").Value = 1144.6364
$ws.Range("K110").Value = 851.2727
$ws.Range("L110").Value = 1144.6364
$ws.Range("M110").Value = 1193.7273
$ws.Range("N110").Value = -5234.6364

$ws.Range("H116").Value = 2102799
$ws.Range("I116").Value = 2332.8572
$ws.Range("K116").Value = 2332.8572
$ws.Range("M116").Value = -38.85719999999992

$ws.Range("H132").Value = 1770.2195
$ws.Range("I132").Value = 1647.5333
$ws.Range("K132").Value = 4942.5999
$ws.Range("M132").Value = -2412.5999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2102799
$ws.Range("I3").Value = 2332.8572
$ws.Range("K3").Value = 2332.8572
$ws.Range("M3").Value = -2218.8572

$ws.Range("H105").Value = 1858.6428
$ws.Range("I105").Value = 1950.8334
$ws.Range("J105").Value = 1305.5
$ws.Range("K105").Value = 1950.8334
$ws.Range("L105").Value = 1305.5
$ws.Range("M105").Value = -203.8334
$ws.Range("N105").Value = -4799.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3275
$ws.Range("I16").Value = 3242.8572
$ws.Range("J16").Value = 3500
$ws.Range("K16").Value = 3242.8572
$ws.Range("L16").Value = 3500
$ws.Range("M16").Value = -2955.8572
$ws.Range("N16").Value = -4074

$ws.Range("H86").Value = 2477.0312
$ws.Range("I86").Value = 2098.1765
$ws.Range("J86").Value = 2906.4
$ws.Range("K86").Value = 2098.1765
$ws.Range("L86").Value = 2906.4
$ws.Range("M86").Value = -975.1765
$ws.Range("N86").Value = -5152.4

$ws.Range("H89").Value = 2477.0312
$ws.Range("I89").Value = 2098.1765
$ws.Range("J89").Value = 2906.4
$ws.Range("K89").Value = 10490.8825
$ws.Range("L89").Value = 14532
$ws.Range("M89").Value = -4874.8825
$ws.Range("N89").Value = -25764

$ws.Range("H99").Value = 1348.4166
$ws.Range("I99").Value = 1148.8334
$ws.Range("J99").Value = 1548
$ws.Range("K99").Value = 1148.8334
$ws.Range("L99").Value = 1548
$ws.Range("M99").Value = 349.1666
$ws.Range("N99").Value = -4544

$ws.Range("H105").Value = 802.75
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 802.75
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 802.75
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -4296.75

$ws.Range("H113").Value = 3275
$ws.Range("I113").Value = 3242.8572
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 3242.8572
$ws.Range("L113").Value = 3500
$ws.Range("M113").Value = -1072.8572
$ws.Range("N113").Value = -7840

$ws.Range("H126").Value = 1348.4166
$ws.Range("I126").Value = 1148.8334
$ws.Range("J126").Value = 1548
$ws.Range("K126").Value = 3446.5002
$ws.Range("L126").Value = 4644
$ws.Range("M126").Value = -976.5001999999999
$ws.Range("N126").Value = -9584

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1981.8182
$ws.Range("I3").Value = 800
$ws.Range("K3").Value = 2400
$ws.Range("M3").Value = -2288

$ws.Range("H113").Value = 537.5319
$ws.Range("I113").Value = 501.82608
$ws.Range("J113").Value = 571.75
$ws.Range("K113").Value = 1505.47824
$ws.Range("L113").Value = 1715.25
$ws.Range("M113").Value = 664.5217600000001
$ws.Range("N113").Value = -6055.25

$ws.Range("H131").Value = 2384091.8
$ws.Range("I131").Value = 4535.8335
$ws.Range("J131").Value = 3335914.2
$ws.Range("K131").Value = 13607.5005
$ws.Range("L131").Value = 10007742.6
$ws.Range("M131").Value = -8567.500499999998
$ws.Range("N131").Value = -10017822.6

$ws.Range("H133").Value = 7183.25
$ws.Range("I133").Value = 3384
$ws.Range("J133").Value = 8183.0527
$ws.Range("K133").Value = 10152
$ws.Range("L133").Value = 24549.1581
$ws.Range("M133").Value = -5092
$ws.Range("N133").Value = -34669.1581

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1849.6875
$ws.Range("I113").Value = 1833.8182
$ws.Range("K113").Value = 1833.8182
$ws.Range("M113").Value = 336.1818000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2443.889
$ws.Range("I40").Value = 2416
$ws.Range("J40").Value = 2499.6667
$ws.Range("K40").Value = 2416
$ws.Range("L40").Value = 2499.6667
$ws.Range("M40").Value = -2280
$ws.Range("N40").Value = -2771.6667

$ws.Range("H61").Value = 1929.1428
$ws.Range("I61").Value = 1156.4445
$ws.Range("J61").Value = 3320
$ws.Range("K61").Value = 1156.4445
$ws.Range("L61").Value = 3320
$ws.Range("M61").Value = -954.4445000000001
$ws.Range("N61").Value = -3724

$ws.Range("H68").Value = 1830.15
$ws.Range("I68").Value = 1792.3077
$ws.Range("J68").Value = 1900.4286
$ws.Range("K68").Value = 1792.3077
$ws.Range("L68").Value = 1900.4286
$ws.Range("M68").Value = -1043.3077
$ws.Range("N68").Value = -3398.4286

$ws.Range("H71").Value = 1830.15
$ws.Range("I71").Value = 1792.3077
$ws.Range("J71").Value = 1900.4286
$ws.Range("K71").Value = 8961.538500000001
$ws.Range("L71").Value = 9502.143
$ws.Range("M71").Value = -5217.538500000001
$ws.Range("N71").Value = -16990.143

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H113").Value = 1929.1428
$ws.Range("I113").Value = 1156.4445
$ws.Range("J113").Value = 3320
$ws.Range("K113").Value = 1156.4445
$ws.Range("L113").Value = 3320
$ws.Range("M113").Value = 1013.5555
$ws.Range("N113").Value = -7660

$ws.Range("H122").Value = 16467.334
$ws.Range("I122").Value = 23803.2
$ws.Range("J122").Value = 7297.5
$ws.Range("K122").Value = 71409.60000000001
$ws.Range("L122").Value = 21892.5
$ws.Range("M122").Value = -68959.60000000001
$ws.Range("N122").Value = -26792.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7259.1
$ws.Range("I62").Value = 8869.200000000001
$ws.Range("J62").Value = 5649
$ws.Range("K62").Value = 8869.200000000001
$ws.Range("L62").Value = 5649
$ws.Range("M62").Value = -8245.200000000001
$ws.Range("N62").Value = -6897

$ws.Range("H65").Value = 7259.1
$ws.Range("I65").Value = 8869.200000000001
$ws.Range("J65").Value = 5649
$ws.Range("K65").Value = 44346
$ws.Range("L65").Value = 28245
$ws.Range("M65").Value = -41226
$ws.Range("N65").Value = -34485

$ws.Range("H100").Value = 486.875
$ws.Range("I100").Value = 478
$ws.Range("J100").Value = 501.66666
$ws.Range("K100").Value = 956
$ws.Range("L100").Value = 1003.33332
$ws.Range("M100").Value = -415
$ws.Range("N100").Value = -2085.33332

$ws.Range("H113").Value = 489.7647
$ws.Range("J113").Value = 433.16666
$ws.Range("L113").Value = 1299.49998
$ws.Range("N113").Value = -5639.499980000001

$ws.Range("H126").Value = 4087.6667
$ws.Range("I126").Value = 2957.1538
$ws.Range("J126").Value = 11436
$ws.Range("K126").Value = 8871.4614
$ws.Range("L126").Value = 34308
$ws.Range("M126").Value = -6401.4614
$ws.Range("N126").Value = -39248

$ws.Range("H136").Value = 18292.166
$ws.Range("I136").Value = 21491.1
$ws.Range("J136").Value = 2297.5
$ws.Range("K136").Value = 64473.3
$ws.Range("L136").Value = 6892.5
$ws.Range("M136").Value = -61923.3
$ws.Range("N136").Value = -11992.5
